# Update crypto symbol list data (price / volume(1h) / coin+link for re-ranked rows)
# as refreshed by the GitHub Actions job on Mon Jan  2 19:23:46 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '246.41'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.45%'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '29.79'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '9.59%'

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.175'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.63%'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05728'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.42%'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.594'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '1.49%'

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8557'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '4.49%'

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8670'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '0.31%'

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1364'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '2.37%'

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '1.95%'

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.02937'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '3.99%'

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09385'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.10%'

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001514'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.96%'

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.04168'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '2.68%'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0005998'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-94.05%'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005986'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-1.86%'

# Row 17
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '5,073.18%'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.489'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.55%'

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.098'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '3.03%'

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.268'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-2.07%'

# Row 21
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-0.31%'

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.03463'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '8.47%'

# Row 23
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '1.05%'

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.481'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-2.21%'

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '0.49%'

# Row 26: 'HotbitToken' -> 'BitKan'
$ws.Range('B26').Value = 'BitKan'
$ws.Range('C26').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.001232'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '1.17%'

# Row 27: 'BitKan' -> 'HotbitToken'
$ws.Range('B27').Value = 'HotbitToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.005020'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '12.03%'

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '22.32%'

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03751'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '0.63%'

# Row 41: 'BKEXToken' -> 'KickToken'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.005739'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-2.80%'

# Row 42: 'CEJI' -> 'BKEXToken'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1071'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '1.17%'

# Row 43: 'KickToken' -> 'CEJI'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002000'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-12.98%'

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.009582'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '2.07%'

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005220'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '1.33%'

# Row 46
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.07%'

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.06468'

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '0.12%'

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.07%'

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.07%'

Write-Host "Symbol list updated with GitHub Actions refresh data"
